# Apply "Horarios actualizados Linea 141 - 582" update to the workbook.
#
# Workbook has 3 sheets:
#   1) LP1912      (main schedule)
#   2) LP1912-215  (sub schedule)
#   3) 6203-6173   (other line, only timestamp changes)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header info
$ws1.Range("A2").Value = "Última actualización: 11:53:44"
$ws1.Range("A3").Value = "Total filas: 147"

# Swap Linea values for two pairs of rows (arrival times tie, order of the
# two stops got swapped by the scraper)
$ws1.Range("C106").Value = "10_OLMOS"
$ws1.Range("C107").Value = "16_SANTA ANA"

$ws1.Range("C133").Value = "14_ABASTO"
$ws1.Range("C134").Value = "27_EL RETIRO"

# Rows 138-152: new snapshot of the tail of the schedule table (4 new rows
# appended/inserted, and re-sorted by arrival time)
$rows1 = @(
    @(138, "11:53:44", "12:36", "23_HERNANDEZ",       43,  "LP1912"),
    @(139, "11:33:52", "12:47", "14_ABASTO",           74,  "LP1912"),
    @(140, "10:49:38", "12:48", "16_SANTA ANA",        119, "LP1912"),
    @(141, "11:33:52", "12:48", "15X38_ABASTO",        75,  "LP1912"),
    @(142, "11:33:52", "13:02", "11_ETCHEVERRY",       89,  "LP1912"),
    @(143, "11:33:52", "13:03", "215C_EL PATO",        90,  "LP1912"),
    @(144, "11:13:15", "13:03", "11_ETCHEVERRY",       110, "LP1912"),
    @(145, "11:46:32", "13:04", "215C_EL PATO",        78,  "LP1912"),
    @(146, "11:33:52", "13:13", "16_SANTA ANA",        100, "LP1912"),
    @(147, "11:33:52", "13:17", "10_OLMOS",            104, "LP1912"),
    @(148, "11:53:44", "13:21", "23_HERNANDEZ",        88,  "LP1912"),
    @(149, "11:33:52", "13:25", "16_P MOR-SANTA ANA",  112, "LP1912"),
    @(150, "11:53:44", "13:32", "215A_EL PATO",        99,  "LP1912"),
    @(151, "11:46:32", "13:33", "215A_EL PATO",        107, "LP1912"),
    @(152, "11:53:44", "13:47", "225_GOMEZ",           114, "LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Range("A$r").Value = $row[1]
    $ws1.Range("B$r").Value = $row[2]
    $ws1.Range("C$r").Value = $row[3]
    $ws1.Range("D$r").Value = $row[4]
    $ws1.Range("E$r").Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 11:53:44"
$ws2.Range("A3").Value = "Total filas: 28"

# Insert a new row before the old last row (old row 32 becomes row 33)
$ws2.Range("A33").Value = "11:46:32"
$ws2.Range("B33").Value = "13:33"
$ws2.Range("C33").Value = "215A_EL PATO"
$ws2.Range("D33").Value = 107
$ws2.Range("E33").Value = "LP1912"

$ws2.Range("A32").Value = "11:53:44"
$ws2.Range("B32").Value = "13:32"
$ws2.Range("C32").Value = "215A_EL PATO"
$ws2.Range("D32").Value = 99
$ws2.Range("E32").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 11:53:44"
